$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3770079612731934
$ws.Range("B1").Value = 0.4749196171760559
$ws.Range("C1").Value = 0.6900337338447571
$ws.Range("D1").Value = 3.03247332572937
$ws.Range("E1").Value = 5.327220439910889
